$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before current row 6 ("Pet tracker Route Design" link row)
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# Row 6: /pets/pettype
$ws.Range("A6").Value = "/pets/pettype"
$ws.Range("B6").Value = "GET"
$ws.Range("C6").Value = "Get all the pet types"
$ws.Range("D6").Value = "id: AUTO ID`npet_type: Cat, Dog, Hedgehog, Oceanic, Reptiles, Exotics, Others"

# Row 7: /pets/createpet
$ws.Range("A7").Value = "/pets/createpet"
$ws.Range("B7").Value = "POST"
$ws.Range("C7").Value = "Create a new pet"
$ws.Range("D7").Value = "petname: STRING,`npetimg: STRING,`npettype_id: AUTO"

# Copy style from row 5 (A:C) for rows 6/7, then set row height
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C7").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").Copy()
$ws.Range("D6:D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("D6:D7").WrapText = $true

$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 43.2

# Column widths
$ws.Columns.Item(1).ColumnWidth = 69.44140625
$ws.Columns.Item(3).ColumnWidth = 39

# View settings
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("F6").Select()
